$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

$xlPasteValues = -4163

# 1. Fix the duplicated word in the title of task Id=18 (row 9 on Active,
#    "have have as little border space..." -> "have as little border space...")
#    before moving the row, so the corrected text carries over.
$active.Cells.Item(9, 2).Value = "have as little border space between elements as possible"

# 2. Grab literal copies of an existing "Done" status cell and an existing
#    "3/3/2018" done-date cell from the Inactive sheet so that pasting them
#    into the new rows keeps them as plain text (matching how every other
#    row on this sheet stores Status/Done), instead of Excel re-parsing a
#    freshly-typed "3/3/2018" string into a date serial. Both scratch cells
#    live on row 1 (in spare columns) so later inserts of rows below row 1
#    never shift them.
$inactive.Range("C2").Copy()
$doneStatusRange = $inactive.Range("Y1")
$doneStatusRange.PasteSpecial($xlPasteValues)

$inactive.Range("F2").Copy()
$doneDateRange = $inactive.Range("Z1")
$doneDateRange.PasteSpecial($xlPasteValues)

# 3. Make room at the top of the Inactive sheet's data for the two tasks
#    that are being marked Done (final order: Id 18, then Id 16). A row
#    Insert() inherits the formatting of the row above it (the bold header
#    row here), so clear that back to the plain/default style used by every
#    other data row on this sheet.
$inactive.Rows.Item(2).Insert()
$inactive.Range("A2:F2").ClearFormats()
$inactive.Rows.Item(2).Insert()
$inactive.Range("A2:F2").ClearFormats()

# 4. Copy each source row's Id/Title/Category/Created (A:E) from Active over
#    to Inactive as plain values (preserves text-vs-number/date typing
#    exactly as stored, with no reinterpretation).
$active.Range("A9:E9").Copy()
$inactive.Range("A2:E2").PasteSpecial($xlPasteValues)

$active.Range("A7:E7").Copy()
$inactive.Range("A3:E3").PasteSpecial($xlPasteValues)

# 5. Fill in Status ("Done") and the Done-date ("3/3/2018") for both new rows.
$doneStatusRange.Copy()
$inactive.Range("C2").PasteSpecial($xlPasteValues)
$doneStatusRange.Copy()
$inactive.Range("C3").PasteSpecial($xlPasteValues)

$doneDateRange.Copy()
$inactive.Range("F2").PasteSpecial($xlPasteValues)
$doneDateRange.Copy()
$inactive.Range("F3").PasteSpecial($xlPasteValues)

$doneStatusRange.ClearContents()
$doneDateRange.ClearContents()

# 6. Remove the two rows from Active now that they live on Inactive
#    (delete the higher row index first so the lower one's position
#    doesn't shift before it is deleted).
$active.Rows.Item(9).Delete()
$active.Rows.Item(7).Delete()
